$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to text format first so numeric-looking strings
# (e.g. "59.203.20", "5.52") are stored as literal text, not coerced to
# Double values with floating-point artifacts. Restore afterward so the
# cell style indices are unaffected (back to default style 0).
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '59.203.20'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '2.324.95'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '541.98'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = '132.11'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.584'
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('D9').Value = '2.320.67'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D11').Value = '5.52'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '0.150'
$ws.Range('D13').Value = '0.333'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '23.85'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '2.738.35'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = '59.206.31'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '2.323.47'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('D19').Value = '10.61'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = '4.18'
$ws.Range('E20').Value = '  -2.07%  '
$ws.Range('D21').Value = '314.19'
$ws.Range('D22').Value = '6.63'
$ws.Range('E22').Value = '  +3.06%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '62.69'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').Value = '0.175'
$ws.Range('E25').Value = '  +3.93%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '8.00'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').Value = '171.33'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('B30').Value = 'SuiNetwork'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D30').Value = '1.18'
$ws.Range('E30').Value = '  +8.55%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.71'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  +3.32%  '
$ws.Range('D33').Value = '5.88'
$ws.Range('E33').Value = '  +3.13%  '
$ws.Range('E34').Value = '  +17.06%  '
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '17.86'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = '4.06'
$ws.Range('E39').Value = '  +3.57%  '
$ws.Range('D40').Value = '317.02'
$ws.Range('E40').Value = '  +9.11%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '38.04'
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.53'
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').Value = '142.21'
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('D44').Value = '3.44'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.559'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0493'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = '18.36'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').Value = '0.0210'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('E51').Value = '  -0.32%  '

# Restore default formatting/style so only cell values differ from the original.
$numRange.NumberFormat = "General"
$numRange.Style = "Normal"

